$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Graduate Students")

$values = @{
    2  = "All students"
    3  = "Male"
    4  = "Female"
    5  = "U.S. citizens and permanent residents"
    6  = "Hispanic or Latino"
    7  = "Not Hispanic or Latino"
    8  = "American Indian or Alaska Native"
    9  = "Asian"
    10 = "Black or African American"
    11 = "Native Hawaiian or Other Pacific Islander"
    12 = "White"
    13 = "More than one race"
    14 = "Unknown ethnicity and race"
    15 = "Temporary visa holders"
    16 = "Science and engineering"
    17 = "Science"
    18 = "Agricultural and veterinary sciences"
    19 = "Biological and biomedical sciences"
    20 = "Communication"
    21 = "Computer and information sciences"
    22 = "Family and consumer sciences and human sciences"
    23 = "Geosciences, atmospheric sciences, and ocean sciences"
    24 = "Mathematics and statistics"
    25 = "Multidisciplinary and interdisciplinary studies"
    26 = "Natural resources and conservation"
    27 = "Psychology"
    28 = "Physical sciences"
    29 = "Social sciences"
    30 = "Engineering"
    31 = "Aerospace, aeronautical, and astronautical engineering"
    32 = "Biological, biomedical, and biosystems engineering"
    33 = "Chemical, petroleum, and chemical-related engineering"
    34 = "Civil, environmental, transportation and related engineering fields"
    35 = "Electrical, electronics, communications and computer engineering"
    36 = "Industrial, manufacturing, systems engineering and operations research"
    37 = "Mechanical engineering"
    38 = "Metallurgical, mining, materials and related engineering fields"
    39 = "Other engineering"
    40 = "Health"
    41 = "Clinical medicine"
    42 = "Other health"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
